$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the style of the "Magnetic Force Between Two Balls" row (old row 73)
# so that after the re-sort it carries the same visual style (index 17 equiv.)
# as its sorted neighbours instead of the stray "thin-gray-border" variant.
$ws.Range("C5").Copy($ws.Range("C73"))

# --- Correct the mis-entered LeetCode number for "Thousand Separator"
$ws.Range("A76").Value = 1556

# --- Add the two new Biweekly Contest 33 problems
$ws.Range("A77").Value = 1557
$ws.Range("B77").Value = "Minimum Number of Vertices to Reach All Nodes"
$ws.Range("C3:D3").Copy($ws.Range("C77:D77"))

$ws.Range("A78").Value = 1558
$ws.Range("B78").Value = "Minimum Number of Function Calls to Make Target Array"
$ws.Range("C4:D4").Copy($ws.Range("C78:D78"))

# --- Re-sort the whole table (A2:D78) by S.No ascending, as the sheet
# always keeps the problem list ordered by LeetCode number.
$rng = $ws.Range("A2:D78")
$key = $ws.Range("A2")
$rng.Sort($key, 1, $null, $null, 1, $null, 1, 1)

# --- Restore the selection that was active when the author saved the file.
$ws.Range("H65").Select()
